$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows right before the existing "Are you receiving any sort of
# public or employment benefits?" row (current row 68), pushing it and
# everything below it down by two rows.
$ws.Rows("68:69").Insert()

# New row 68: partner's employment status PRIOR to the pandemic.
$ws.Range("A68").Value = "What was your partner's employment status prior to the coronavirus (COVID-19) Pandemic? If you do not have a partner, please select not applicable. "
$ws.Range("B68").Value = "• Working`n• Unemployed or laid off`n• Temporarily out of work or furloughed`n• Other`n• Not applicable"
$ws.Range("C68").Value = "Income and Employment"
$ws.Range("D68").Value = "RAPID Team Modified from U.S. Census "
$ws.Rows("68").RowHeight = 80

# New row 69: partner's CURRENT employment status.
$ws.Range("A69").Value = "What is your partner's current employment status? If you do not have a partner, please select not applicable. "
$ws.Range("B69").Value = "• Working`n• Unemployed or laid off`n• Temporarily out of work or furloughed`n• Other`n• Not applicable"
$ws.Range("C69").Value = "Income and Employment"
$ws.Range("D69").Value = "RAPID Team Modified from U.S. Census "
$ws.Rows("69").RowHeight = 80

# Match the author's final selection/view state.
$ws.Range("C69").Select()
